# "Proyectos Automatización" — simplify header row:
#   - A1/C1 used to share the same "Proyecto" label; they become distinct
#     "Proyecto ID" (with a leading line break, matching the source sheet)
#     and "Proyecto Descripción" headers.
#   - The bold/shaded header style is dropped in favour of the sheet's plain
#     default style, and the header row's custom height is cleared.
#   - The active selection moves from I17 to A4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New header text for columns A and C (B, D, E headers are unchanged).
$ws.Range("A1").Value = "`nProyecto ID"
$ws.Range("C1").Value = "Proyecto Descripción"

# 2) Strip the bold font + fill/alignment header styling back to the
#    workbook's plain "Normal" style for the whole header row.
$ws.Range("A1:E1").Style = "Normal"

# 3) Drop the header row's custom 15.6pt height, reverting to the sheet's
#    default row height.
$ws.Rows.Item(1).AutoFit()

# 4) Move the selection to A4 (was I17).
$null = $ws.Range("A4").Select()
